$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.958.36'
$ws.Range('E2').Value = '  +1.56%  '
$ws.Range('D3').Value = '3.776.30'
$ws.Range('E3').Value = '  -0.21%  '
$ws.Range('E4').Value = '  -0.17%  '
$ws.Range('D5').Value = "'628.47"
$ws.Range('E5').Value = '  +4.16%  '
$ws.Range('D6').Value = "'165.14"
$ws.Range('E6').Value = '  +1.04%  '
$ws.Range('D7').Value = '3.772.58'
$ws.Range('E7').Value = '  -0.23%  '
$ws.Range('E8').Value = '  +0.12%  '
$ws.Range('E9').Value = '  +1.36%  '
$ws.Range('D10').Value = "'0.159"
$ws.Range('E10').Value = '  +0.96%  '
$ws.Range('E11').Value = '  +2.51%  '
$ws.Range('D12').Value = "'6.81"
$ws.Range('E12').Value = '  +0.07%  '
$ws.Range('D13').Value = "'0.0000244"
$ws.Range('E13').Value = '  -0.80%  '
$ws.Range('D14').Value = "'35.19"
$ws.Range('E14').Value = '  +0.38%  '
$ws.Range('D15').Value = '4.411.63'
$ws.Range('E15').Value = '  -0.17%  '
$ws.Range('D16').Value = '3.773.31'
$ws.Range('E16').Value = '  -0.50%  '
$ws.Range('D17').Value = '68.986.53'
$ws.Range('E17').Value = '  +1.63%  '
$ws.Range('D18').Value = "'17.58"
$ws.Range('E18').Value = '  -3.30%  '
$ws.Range('E19').Value = '  -1.12%  '
$ws.Range('D20').Value = "'7.03"
$ws.Range('E20').Value = '  +0.41%  '
$ws.Range('E21').Value = '  +1.81%  '
$ws.Range('E22').Value = '  +0.78%  '
$ws.Range('E23').Value = '  +2.39%  '
$ws.Range('D24').Value = "'82.98"
$ws.Range('E24').Value = '  -0.23%  '
$ws.Range('E25').Value = '  +0.24%  '
$ws.Range('D26').Value = "'12.00"
$ws.Range('E26').Value = '  +1.14%  '
$ws.Range('E27').Value = '  +3.42%  '
$ws.Range('D28').Value = "'10.02"
$ws.Range('E28').Value = '  +1.28%  '
$ws.Range('E29').Value = '  +0.01%  '
$ws.Range('D30').Value = '3.925.84'
$ws.Range('E30').Value = '  -0.21%  '
$ws.Range('E31').Value = '  +2.74%  '
$ws.Range('D32').Value = "'2.23"
$ws.Range('E32').Value = '  +2.31%  '
$ws.Range('D33').Value = "'7.12"
$ws.Range('E33').Value = '  -1.16%  '
$ws.Range('D34').Value = "'28.72"
$ws.Range('E34').Value = '  -0.89%  '
$ws.Range('D35').Value = "'0.172"
$ws.Range('E35').Value = '  +16.47%  '
$ws.Range('E36').Value = '  +0.06%  '
$ws.Range('D37').Value = '3.726.93'
$ws.Range('E37').Value = '  -0.22%  '
$ws.Range('E38').Value = '  +0.29%  '
$ws.Range('E39').Value = '  +2.08%  '
$ws.Range('D40').Value = "'3.29"
$ws.Range('E40').Value = '  +2.41%  '
$ws.Range('E41').Value = '  +0.14%  '
$ws.Range('D42').Value = "'0.966"
$ws.Range('E42').Value = '  -1.12%  '
$ws.Range('E43').Value = '  +0.00%  '
$ws.Range('E44').Value = '  -0.14%  '
$ws.Range('D45').Value = "'155.85"
$ws.Range('E45').Value = '  +2.15%  '
$ws.Range('D46').Value = "'42.91"
$ws.Range('E46').Value = '  -1.76%  '
$ws.Range('E47').Value = '  +0.46%  '
$ws.Range('D48').Value = "'46.63"
$ws.Range('E48').Value = '  -0.94%  '
$ws.Range('E49').Value = '  +3.78%  '
$ws.Range('D50').Value = "'8.36"
$ws.Range('E50').Value = '  +0.79%  '
$ws.Range('D51').Value = "'1.36"
$ws.Range('E51').Value = '  -1.14%  '
